$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new log row mirrors the previous row's layout/formatting (SKIPPED,
# same message/URL, blank Saved PDF, 0 rows appended, blank Total Rows),
# so copy row 61 (values + styles) down to row 62, then patch only the
# two timestamp cells that actually changed.
$src = $ws.Range("A61:H61")
$dst = $ws.Range("A62:H62")
$src.Copy($dst)

$ws.Cells.Item(62, 1).Value = "2025-08-27 03:47:21 UTC"
$ws.Cells.Item(62, 2).Value = "2025-08-27 09:17:21 IST"
